$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 86) with the same shape as the existing rows.
$row = 86

# Column A holds a date written as plain text (e.g. "2025/10/09"). Force
# the cell to Text format before assigning it so Excel does not silently
# convert it into a date serial number, then clear the formatting again
# so the cell ends up with no explicit style - matching the rest of the
# data rows in the sheet.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/09"
$ws.Cells.Item($row, 1).ClearFormats()

# Column B is the weekday (kanji) as plain text.
$ws.Cells.Item($row, 2).Value = "木"

# Columns C and D are numeric values.
$ws.Cells.Item($row, 3).Value = 20
$ws.Cells.Item($row, 4).Value = 201
